$d = $word.ActiveDocument

# "Compartición de datos" -> "Compartir tus datos" (bold runs stay bold) in
# the "Paragraph 4" heading. Both find/replace calls below are scoped to a
# fresh, unique phrase-match range so the very common Spanish word " de "
# elsewhere in the document is left untouched.

# Step 1: "Compartición" -> "Compartir"
$r1 = $d.Content
$r1.Find.Execute("Compartición de datos", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$scope1 = $d.Range($r1.Start, $r1.End)
$scope1.Find.Execute("Compartición", $false, $false, $false, $false, $false, `
                      $true, 1, $false, "Compartir", 1)

# Step 2: " de " -> " tus " (re-locate the now-updated unique phrase first,
# since the text shrank by three characters after step 1).
$r2 = $d.Content
$r2.Find.Execute("Compartir de datos", $false, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$scope2 = $d.Range($r2.Start, $r2.End)
$scope2.Find.Execute(" de ", $false, $false, $false, $false, $false, `
                      $true, 1, $false, " tus ", 1)
